$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.028.73'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '1.643.93'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +1.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.07%  '
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('E7').Value = '  +1.01%  '
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.64'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('E11').Value = '  +1.01%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.698.99'
$ws.Range('E12').Value = '  +4.41%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.30'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.58%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.872.24'
$ws.Range('E14').Value = '  +0.93%  '
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').Value = '0.0₃0765'
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Value = '26.090.03'
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '193.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.42'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.66%  '
$ws.Range('B26').Value = 'BinanceUSD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.25%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.131'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.91'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('E33').Value = '  +1.23%  '
$ws.Range('E35').Value = '  -3.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.906'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.64%  '
$ws.Range('D37').Value = '1.133.75'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.797'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.53%  '
$ws.Range('D44').Value = '1.781.19'
$ws.Range('E44').Value = '  +0.91%  '
$ws.Range('D45').Value = '0.0₆0117'
$ws.Range('E45').Value = '  +5.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.58'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('E47').Value = '  +0.98%  '
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.70'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.416'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.66%  '
$ws.Range('E51').Value = '  -0.31%  '
